$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 88.72291666666666
$ws.Range("H2").Value = 266.16875
$ws.Range("I2").Value = 0.7675060578750151
$ws.Range("J2").Value = 0.7675060578750152
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 27.592778
$ws.Range("N2").Value = 82.778334
$ws.Range("O2").Value = 0.2684079248986126
$ws.Range("P2").Value = 0.2684079248986126
$ws.Range("Q2").Value = 2448.111743095833
$ws.Range("R2").Value = 22033.0056878625
$ws.Range("S2").Value = 0.2060047083413473
$ws.Range("T2").Value = 0.2060047083413473

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 88.72291666666666
$ws.Range("H3").Value = 266.16875
$ws.Range("I3").Value = 0.7675060578750151
$ws.Range("J3").Value = 0.7675060578750152
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 65.63594833333333
$ws.Range("N3").Value = 196.907845
$ws.Range("O3").Value = 0.6384717294830753
$ws.Range("P3").Value = 0.6384717294830752
$ws.Range("Q3").Value = 5823.412774315972
$ws.Range("R3").Value = 52410.71496884375
$ws.Range("S3").Value = 0.4900309201601982
$ws.Range("T3").Value = 0.4900309201601982

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 88.72291666666666
$ws.Range("H4").Value = 266.16875
$ws.Range("I4").Value = 0.7675060578750151
$ws.Range("J4").Value = 0.7675060578750152
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 9.572925333333334
$ws.Range("N4").Value = 28.718776
$ws.Range("O4").Value = 0.09312034561831214
$ws.Range("P4").Value = 0.09312034561831213
$ws.Range("Q4").Value = 849.3378566055555
$ws.Range("R4").Value = 7644.040709450001
$ws.Range("S4").Value = 0.07147042937346969
$ws.Range("T4").Value = 0.07147042937346969

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 17.91585
$ws.Range("H5").Value = 53.74755
$ws.Range("I5").Value = 0.1549827702197958
$ws.Range("J5").Value = 0.1549827702197958
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 27.592778
$ws.Range("N5").Value = 82.778334
$ws.Range("O5").Value = 0.2684079248986126
$ws.Range("P5").Value = 0.2684079248986126
$ws.Range("Q5").Value = 494.3480717313001
$ws.Range("R5").Value = 4449.132645581701
$ws.Range("S5").Value = 0.04159860374973388
$ws.Range("T5").Value = 0.04159860374973389

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 17.91585
$ws.Range("H6").Value = 53.74755
$ws.Range("I6").Value = 0.1549827702197958
$ws.Range("J6").Value = 0.1549827702197958
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 65.63594833333333
$ws.Range("N6").Value = 196.907845
$ws.Range("O6").Value = 0.6384717294830753
$ws.Range("P6").Value = 0.6384717294830752
$ws.Range("Q6").Value = 1175.92380494775
$ws.Range("R6").Value = 10583.31424452975
$ws.Range("S6").Value = 0.09895211734231109
$ws.Range("T6").Value = 0.09895211734231109

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 17.91585
$ws.Range("H7").Value = 53.74755
$ws.Range("I7").Value = 0.1549827702197958
$ws.Range("J7").Value = 0.1549827702197958
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 9.572925333333334
$ws.Range("N7").Value = 28.718776
$ws.Range("O7").Value = 0.09312034561831214
$ws.Range("P7").Value = 0.09312034561831213
$ws.Range("Q7").Value = 171.5070943332
$ws.Range("R7").Value = 1543.5638489988
$ws.Range("S7").Value = 0.01443204912775084
$ws.Range("T7").Value = 0.01443204912775084

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 8.960212333333333
$ws.Range("H8").Value = 26.880637
$ws.Range("I8").Value = 0.077511171905189
$ws.Range("J8").Value = 0.07751117190518901
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 27.592778
$ws.Range("N8").Value = 82.778334
$ws.Range("O8").Value = 0.2684079248986126
$ws.Range("P8").Value = 0.2684079248986126
$ws.Range("Q8").Value = 247.2371497465286
$ws.Range("R8").Value = 2225.134347718758
$ws.Range("S8").Value = 0.02080461280753142
$ws.Range("T8").Value = 0.02080461280753142

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 8.960212333333333
$ws.Range("H9").Value = 26.880637
$ws.Range("I9").Value = 0.077511171905189
$ws.Range("J9").Value = 0.07751117190518901
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 65.63594833333333
$ws.Range("N9").Value = 196.907845
$ws.Range("O9").Value = 0.6384717294830753
$ws.Range("P9").Value = 0.6384717294830752
$ws.Range("Q9").Value = 588.1120337663627
$ws.Range("R9").Value = 5293.008303897265
$ws.Range("S9").Value = 0.04948869198056598
$ws.Range("T9").Value = 0.04948869198056598

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 8.960212333333333
$ws.Range("H10").Value = 26.880637
$ws.Range("I10").Value = 0.077511171905189
$ws.Range("J10").Value = 0.07751117190518901
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 9.572925333333334
$ws.Range("N10").Value = 28.718776
$ws.Range("O10").Value = 0.09312034561831214
$ws.Range("P10").Value = 0.09312034561831213
$ws.Range("Q10").Value = 85.77544363781244
$ws.Range("R10").Value = 771.978992740312
$ws.Range("S10").Value = 0.007217867117091606
$ws.Range("T10").Value = 0.007217867117091606
